$d = $word.ActiveDocument

$replacements = @(
    @{ old = "<id>p049r_a1</id>"; new = "<id>p049r_1</id>" },
    @{ old = "<id>p049r_a2</id>"; new = "<id>p049r_2</id>" },
    @{ old = "<id>p049r_a3</id>"; new = "<id>p049r_3</id>" },
    @{ old = "<id>p049r_a4</id>"; new = "<id>p049r_4</id>" }
)

foreach ($rep in $replacements) {
    $rng = $d.Content
    $found = $rng.Find.Execute($rep.old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $rep.new
    }
}
